$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.866
$ws.Range("B4").Value = 7.043000000000001

$ws.Range("B5").Value = 6.38

$ws.Range("A6").Value = -21.14
$ws.Range("B6").Value = 6.816

$ws.Range("A7").Value = -21.089

$ws.Range("A8").Value = -21.046
$ws.Range("B8").Value = 6.161

$ws.Range("A16").Value = -20.998
$ws.Range("B16").Value = 6.653999999999999

$ws.Range("A20").Value = -21.86

$ws.Range("A21").Value = -21.14

$ws.Range("B22").Value = 6.626
